$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 28: new cancellation record (CVR 45357716) ---
# Column A holds a CVR number that LOOKS like a plain integer; format the
# cell as Text first so Excel stores it as a shared string (matching the
# other CVR cells in column A) instead of silently coercing it to a number.
$ws.Cells.Item(28, 1).NumberFormat = "@"
$ws.Cells.Item(28, 1).Value = "45357716"
$ws.Cells.Item(28, 1).Style = "Normal"

$ws.Cells.Item(28, 2).Value = 2023
$ws.Cells.Item(28, 3).Value = 60860
$ws.Cells.Item(28, 4).Value = "EasyCruit"

$ws.Cells.Item(28, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(28, 5).Value = 45168

$ws.Cells.Item(28, 8).Value = "2023Q3"
$ws.Cells.Item(28, 9).Value = "60000-80000"

# --- Row 29: new cancellation record (CVR 33738811) ---
$ws.Cells.Item(29, 1).NumberFormat = "@"
$ws.Cells.Item(29, 1).Value = "33738811"
$ws.Cells.Item(29, 1).Style = "Normal"

$ws.Cells.Item(29, 2).Value = 2023
$ws.Cells.Item(29, 3).Value = 65688
$ws.Cells.Item(29, 4).Value = "Visma Løn & HR"

$ws.Cells.Item(29, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(29, 5).Value = 45211

$ws.Cells.Item(29, 8).Value = "2023Q4"
$ws.Cells.Item(29, 9).Value = "60000-80000"
